$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 6525.25
$ws.Range("I43").Value = 7699
$ws.Range("J43").Value = 6134
$ws.Range("K43").Value = 7699
$ws.Range("L43").Value = 6134
$ws.Range("M43").Value = -7630
$ws.Range("N43").Value = -6272

$ws.Range("H55").Value = 95.5
$ws.Range("J55").Value = 90
$ws.Range("L55").Value = 90
$ws.Range("N55").Value = -518

$ws.Range("H107").Value = 820
$ws.Range("I107").Value = 820
$ws.Range("K107").Value = 820
$ws.Range("M107").Value = 1100

$ws.Range("H137").Value = 5683.3335
$ws.Range("J137").Value = 5927.273
$ws.Range("L137").Value = 17781.819
$ws.Range("N137").Value = -22881.819

$ws.Range("H138").Value = 4401.6333
$ws.Range("I138").Value = 3123.6
$ws.Range("J138").Value = 4657.24
$ws.Range("K138").Value = 9370.799999999999
$ws.Range("L138").Value = 13971.72
$ws.Range("M138").Value = -4230.799999999999
$ws.Range("N138").Value = -24251.72

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 815.1667
$ws.Range("I2").Value = 799.44446
$ws.Range("K2").Value = 799.44446
$ws.Range("M2").Value = -686.44446

$ws.Range("H32").Value = 2831.0293
$ws.Range("I32").Value = 1914.25
$ws.Range("K32").Value = 1914.25
$ws.Range("M32").Value = -1627.25

$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("M74").ClearContents()
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("M77").ClearContents()
$ws.Range("N77").ClearContents()

$ws.Range("H116").Value = 815.1667
$ws.Range("I116").Value = 799.44446
$ws.Range("K116").Value = 799.44446
$ws.Range("M116").Value = 1494.55554

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 815.1667
$ws.Range("I3").Value = 799.44446
$ws.Range("K3").Value = 799.44446
$ws.Range("M3").Value = -685.44446

$ws.Range("H86").Value = 999.6667
$ws.Range("I86").Value = 1087.25
$ws.Range("J86").Value = 824.5
$ws.Range("K86").Value = 1087.25
$ws.Range("L86").Value = 824.5
$ws.Range("M86").Value = 35.75
$ws.Range("N86").Value = -3070.5

$ws.Range("H89").Value = 999.6667
$ws.Range("I89").Value = 1087.25
$ws.Range("J89").Value = 824.5
$ws.Range("K89").Value = 5436.25
$ws.Range("L89").Value = 4122.5
$ws.Range("M89").Value = 179.75
$ws.Range("N89").Value = -15354.5

$ws.Range("H94").Value = 863.6667
$ws.Range("I94").Value = 786
$ws.Range("J94").Value = 941.3333
$ws.Range("K94").Value = 786
$ws.Range("L94").Value = 941.3333
$ws.Range("M94").Value = -335
$ws.Range("N94").Value = -1843.3333

$ws.Range("H105").Value = 2532.6667
$ws.Range("I105").Value = 2371.25
$ws.Range("K105").Value = 2371.25
$ws.Range("M105").Value = -624.25

$ws.Range("H134").Value = 2901.7778
$ws.Range("I134").Value = 2827
$ws.Range("J134").Value = 3500
$ws.Range("K134").Value = 8481
$ws.Range("L134").Value = 10500
$ws.Range("M134").Value = -5946
$ws.Range("N134").Value = -15570

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2545.7334
$ws.Range("I31").Value = 2598.8333
$ws.Range("J31").Value = 2333.3333
$ws.Range("K31").Value = 2598.8333
$ws.Range("L31").Value = 2333.3333
$ws.Range("M31").Value = -2303.8333
$ws.Range("N31").Value = -2923.3333

$ws.Range("H34").Value = 2545.7334
$ws.Range("I34").Value = 2598.8333
$ws.Range("J34").Value = 2333.3333
$ws.Range("K34").Value = 2598.8333
$ws.Range("L34").Value = 2333.3333
$ws.Range("M34").Value = -2396.8333
$ws.Range("N34").Value = -2737.3333

$ws.Range("H86").Value = 73000
$ws.Range("I86").Value = 13333.333
$ws.Range("J86").Value = 252000
$ws.Range("K86").Value = 13333.333
$ws.Range("L86").Value = 252000
$ws.Range("M86").Value = -12210.333
$ws.Range("N86").Value = -254246

$ws.Range("H89").Value = 73000
$ws.Range("I89").Value = 13333.333
$ws.Range("J89").Value = 252000
$ws.Range("K89").Value = 66666.66500000001
$ws.Range("L89").Value = 1260000
$ws.Range("M89").Value = -61050.66500000001
$ws.Range("N89").Value = -1271232

$ws.Range("H99").Value = 1293.8182
$ws.Range("J99").Value = 1037
$ws.Range("L99").Value = 1037
$ws.Range("N99").Value = -4033

$ws.Range("H126").Value = 1293.8182
$ws.Range("J126").Value = 1037
$ws.Range("L126").Value = 3111
$ws.Range("N126").Value = -8051

$ws.Range("H132").Value = 4457.727
$ws.Range("I132").Value = 3879.875
$ws.Range("K132").Value = 11639.625
$ws.Range("M132").Value = -9109.625

$ws.Range("H134").Value = 3983
$ws.Range("I134").Value = 3983
$ws.Range("K134").Value = 11949
$ws.Range("M134").Value = -9414

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1250
$ws.Range("J23").Value = 1250
$ws.Range("L23").Value = 3750
$ws.Range("N23").Value = -4220

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H109").Value = 1242
$ws.Range("I109").Value = 784.75
$ws.Range("J109").Value = 2156.5
$ws.Range("K109").Value = 2354.25
$ws.Range("L109").Value = 6469.5
$ws.Range("M109").Value = -1314.25
$ws.Range("N109").Value = -8549.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3363.9412
$ws.Range("I80").Value = 2182
$ws.Range("J80").Value = 6200.6
$ws.Range("K80").Value = 2182
$ws.Range("L80").Value = 6200.6
$ws.Range("M80").Value = -1184
$ws.Range("N80").Value = -8196.6

$ws.Range("H83").Value = 3363.9412
$ws.Range("I83").Value = 2182
$ws.Range("J83").Value = 6200.6
$ws.Range("K83").Value = 10910
$ws.Range("L83").Value = 31003
$ws.Range("M83").Value = -5918
$ws.Range("N83").Value = -40987

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3485.1428
$ws.Range("I7").Value = 3079.4
$ws.Range("K7").Value = 3079.4
$ws.Range("M7").Value = -2967.4

$ws.Range("H126").Value = 3485.1428
$ws.Range("I126").Value = 3079.4
$ws.Range("K126").Value = 9238.200000000001
$ws.Range("M126").Value = -6768.200000000001

$ws.Range("H136").Value = 49999
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2360.125
$ws.Range("I122").Value = 2126
$ws.Range("K122").Value = 6378
$ws.Range("M122").Value = -3928

$ws.Range("H126").Value = 1197.2778
$ws.Range("I126").Value = 1060.8
$ws.Range("J126").Value = 1879.6666
$ws.Range("K126").Value = 3182.4
$ws.Range("L126").Value = 5638.9998
$ws.Range("M126").Value = -712.3999999999996
$ws.Range("N126").Value = -10578.9998

$ws.Range("H133").Value = 119900
$ws.Range("J133").Value = 119900
$ws.Range("L133").Value = 119900
$ws.Range("N133").Value = -130020

$ws.Range("H136").Value = 12824
$ws.Range("I136").Value = 12824
$ws.Range("K136").Value = 38472
$ws.Range("M136").Value = -35922

Write-Output "Applied all changes"
